# This script updates the 'Price' (column D) and 'Volume(1h)' (column E)
# values in the cryptocurrency table on Sheet1, reflecting refreshed
# market data pulled by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> updated column values (only columns that changed
# are present for each row).
$updates = @{
    2 = @{ D='26.089.65'; E='  -1.01%  ' }
    3 = @{ D='1.666.31'; E='  -1.26%  ' }
    4 = @{ E='  -0.54%  ' }
    5 = @{ D='209.12'; E='  -4.28%  ' }
    6 = @{ D='0.5189'; E='  -5.30%  ' }
    7 = @{ D='1.005'; E='  -0.54%  ' }
    8 = @{ D='0.2642'; E='  -3.16%  ' }
    9 = @{ D='0.06193'; E='  -4.06%  ' }
    10 = @{ D='20.94'; E='  -5.02%  ' }
    11 = @{ D='0.07494'; E='  -2.36%  ' }
    12 = @{ D='1.641.41'; E='  -2.80%  ' }
    13 = @{ D='4.433'; E='  -2.25%  ' }
    14 = @{ D='0.5590'; E='  -3.89%  ' }
    15 = @{ D='0.000007934'; E='  -4.87%  ' }
    16 = @{ D='65.54'; E='  +0.71%  ' }
    17 = @{ D='26.120.30'; E='  -1.15%  ' }
    18 = @{ E='  -0.47%  ' }
    19 = @{ D='4.795'; E='  -3.02%  ' }
    20 = @{ D='10.40'; E='  -5.25%  ' }
    21 = @{ D='185.92'; E='  -2.45%  ' }
    22 = @{ D='6.161'; E='  -1.06%  ' }
    23 = @{ D='1.006'; E='  -0.40%  ' }
    24 = @{ D='146.33'; E='  -2.02%  ' }
    25 = @{ D='0.1252'; E='  -4.68%  ' }
    26 = @{ D='7.553'; E='  -4.21%  ' }
    27 = @{ D='15.71'; E='  +0.10%  ' }
    28 = @{ D='0.06356'; E='  -0.18%  ' }
    29 = @{ D='1.344'; E='  -4.95%  ' }
    30 = @{ D='1.270'; E='  -4.34%  ' }
    31 = @{ D='3.484'; E='  -2.47%  ' }
    32 = @{ D='3.434'; E='  -3.96%  ' }
    33 = @{ D='1.624'; E='  -3.02%  ' }
    34 = @{ D='0.9950'; E='  -4.53%  ' }
    35 = @{ D='2.413'; E='  +0.11%  ' }
    36 = @{ D='0.6024'; E='  -2.38%  ' }
    37 = @{ D='2.705'; E='  -0.57%  ' }
    38 = @{ D='6.084'; E='  -2.41%  ' }
    39 = @{ D='0.01610'; E='  -1.09%  ' }
    40 = @{ D='1.093.62'; E='  -1.92%  ' }
    41 = @{ D='0.8615'; E='  -1.97%  ' }
    42 = @{ E='  -1.05%  ' }
    43 = @{ D='99.50'; E='  -1.42%  ' }
    44 = @{ D='1.819.88'; E='  -1.08%  ' }
    45 = @{ D='0.00000000109'; E='  -0.64%  ' }
    46 = @{ D='56.15'; E='  -1.92%  ' }
    47 = @{ D='1.011'; E='  +0.13%  ' }
    48 = @{ D='0.05248' }
    49 = @{ D='7.922'; E='  -3.26%  ' }
    50 = @{ D='0.4265'; E='  -0.91%  ' }
    51 = @{ D='5.897'; E='  -2.36%  ' }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    if ($vals.ContainsKey('D')) {
        $cellD = $ws.Cells.Item($row, 4)
        # Force text format so numeric-looking strings (e.g. '1.005')
        # are preserved verbatim instead of being parsed as numbers.
        $cellD.NumberFormat = '@'
        $cellD.Value = $vals['D']
    }

    if ($vals.ContainsKey('E')) {
        $cellE = $ws.Cells.Item($row, 5)
        $cellE.NumberFormat = '@'
        $cellE.Value = $vals['E']
    }
}
